# Journal workbook update:
# Add a new day's entries (rows 52-55) describing work done after the
# class computers were changed: re-downloading tooling, fixing small
# issues, and starting the base terrain implementation (Issue #9 still
# "Not started"). Row 53 is also given extra height to fit its longer
# task description.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# --- Row 52: "Setup" day header entry ---------------------------------
$ws.Range("D52").Value = "Downloading VS and Github Desktop again"
$ws.Range("G52").Value = "The class computers got changed"
$ws.Range("A52").Value = "Setup"
$ws.Range("B52").Value = 45957
$ws.Range("C52").Value = 20
$ws.Range("E52").Value = "Finished"
$ws.Range("F52").Value = 0.45833333333333331

# --- Row 55: typed first (matches original authoring order) -----------
$ws.Range("D55").Value = "Implementation of Issue #9 to the program"
$ws.Range("A55").Value = "Coding"
$ws.Range("E55").Value = "Not started"

# --- Row 53: fixes / tweaks --------------------------------------------
$ws.Range("D53").Value = "Fixed small issues and tweaked with some variables"
$ws.Range("A53").Value = "Coding"
$ws.Range("C53").Value = 50
$ws.Range("E53").Value = "Finished"
$ws.Range("F53").Value = 0.49305555555555558

# --- Row 54: base terrain implementation -------------------------------
$ws.Range("D54").Value = "Implementation of the base terrain"
$ws.Range("G54").Value = "30m+, started again 13:10"
$ws.Range("A54").Value = "Coding"
$ws.Range("E54").Value = "In the work"

# Row 53 needed extra height once its task text got longer
$ws.Rows.Item(53).RowHeight = 30
